# Insert a new data row before the current row 46 (date serial 44957,
# Primera, 25, 30000/30000/30000, "$/caja 18 kilos", 1667, 18), pushing it
# and every row below it (old rows 46-90) down by one, to old row 91.
# The freshly inserted row 46 gets a brand new observation (date serial
# 45036, Primera, 140, 22000/22000/22000, "$/caja 16 kilos", 1375, 16).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 46:90 down to 47:91 and make room for the new row 46.
$ws.Rows("46:46").Insert()

# Populate the new row 46 with the new observation.
$ws.Cells.Item(46, 1).Value  = 10
$ws.Cells.Item(46, 2).Value  = 'Vega Modelo de Temuco'
$ws.Cells.Item(46, 3).Value  = 'La Araucanía'
$ws.Cells.Item(46, 4).Value  = 45036
$ws.Cells.Item(46, 5).Value  = 9
$ws.Cells.Item(46, 6).Value  = 'Fruta'
$ws.Cells.Item(46, 7).Value  = 100107
$ws.Cells.Item(46, 8).Value  = 'Otros'
$ws.Cells.Item(46, 9).Value  = 100107011
$ws.Cells.Item(46, 10).Value = 'Tuna'
$ws.Cells.Item(46, 11).Value = 'Sin especificar'
$ws.Cells.Item(46, 12).Value = 'Primera'
$ws.Cells.Item(46, 13).Value = 140
$ws.Cells.Item(46, 14).Value = 22000
$ws.Cells.Item(46, 15).Value = 22000
$ws.Cells.Item(46, 16).Value = 22000
$ws.Cells.Item(46, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(46, 18).Value = 'Provincia de Los Andes'
$ws.Cells.Item(46, 19).Value = 1375
$ws.Cells.Item(46, 20).Value = 16
